$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Vlans value for the Trunk port row from "100:200" to "100-200"
$ws.Range("C3").Value = "100-200"

# Reflect the new active selection on that cell (matches the saved view state)
$ws.Range("C3").Select()
